$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 - Индикатор: updated indicator title/wording (8.10.1 -> 8.10.1.)
$ws.Range("B4").Value = "8.10.1. Число а) филиалов коммерческих банков на 100 000 взрослых и b) банкоматов на 100 000 взрослых "

# B6 - Организация: now also mentions the specific NSK KR department
$ws.Range("B6").Value = "НБ КР, НСК КР (Отдел демографической статистики)"
$ws.Range("B6").Font.Name = "Calibri"

# B9 - Телефон контактного лица: phone list updated (Абдуллаева -> Мааткулова)
$ws.Range("B9").Value = "Сулайманбекова Ж.С.: +996-0312-66-92-56;" + [char]10 + "Лелевкина Э.В.: +996-312-61-14-34;" + [char]10 + "Мааткулова Ж.Б.: +996-312-32-55-46"
$ws.Range("B9").Font.Name = "Calibri"

# B7 - Контактное лицо (лица) / Координатор: contact list updated (Абдуллаева -> Мааткулова)
$ws.Range("B7").Value = "Сулайманбекова Ж.С., начальник Управления методологии надзора и лицензирования банков НБ КР, " + [char]10 + "Лелевкина Э.В. , начальник Управления платежных систем НБ КР;" + [char]10 + "Мааткулова Ж.Б. , заведующая отделом демографической статистики НСК КР."
$ws.Range("B7").Font.Name = "Calibri"

# B10 - Сайт организации: НСК КР site URL changed from stat.kg to stat.gov.kg
$ws.Range("B10").Value = "НБ КР: www.nbkr.kg;" + [char]10 + "НСК КР: www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

# Selection moved from B2 to B10
$ws.Range("B10").Select()
